# Update Excel file for new format:
#  - Rename the existing sheet (matrix data) to "ConflictMatrix".
#  - Insert a brand-new "Setup" sheet in front of it containing the
#    signal-group id / conflict-time lookup table (A1:B10).

$wb = $excel.ActiveWorkbook

# --- Existing sheet becomes "ConflictMatrix" -------------------------------
$matrix = $wb.ActiveSheet
$matrix.Name = "ConflictMatrix"

# --- New "Setup" sheet, inserted before the matrix sheet -------------------
$setup = $wb.Worksheets.Add()
$setup.Name = "Setup"

$ids = @(1.1, 2.1, 5.1, 6.1, 7.1, 8.1, 9.1, 10.1, 11.1, 12.1)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 1
    $setup.Cells.Item($row, 1).Value = $ids[$i]
    $setup.Cells.Item($row, 2).Value = 1.5
}

# Match the highlighted/bordered "Comma" look used for the id column on the
# ConflictMatrix sheet (yellow fill, thin border, 1-decimal comma format).
$idCol = $setup.Range("A1:A10")
$idCol.NumberFormat = "_-* #,##0.0_-;\-* #,##0.0_-;_-* ""-""??_-;_-@_-"
$idCol.Interior.Color = 10284031
$idCol.Borders.LineStyle = 1
$idCol.Borders.Weight = 2
$idCol.Borders.ColorIndex = 64

$setup.Activate()
